# Incorporate addenda a, b, and c.
#
# 1. RS0003 sheet: bump schema_version from 1.0.0 to 2.0.0
# 2. RS0003 sheet: data validations on C5/C16/C19/C27/C28 no longer show
#    input/error messages (InCellDropdown stays on, ShowInput/ShowError off)
# 3. performance_map sheet: new lookup_variables column E = operation_state,
#    a "-" unit row, a comment on E3, and "NORMAL" for every data row (5-144)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. RS0003!C6 schema_version bump
# ---------------------------------------------------------------------------
$rs0003 = $wb.Worksheets.Item("RS0003")
$rs0003.Range("C6").Value = "2.0.0"

# ---------------------------------------------------------------------------
# 2. RS0003 data validations: stop showing input/error messages
# ---------------------------------------------------------------------------
$dvCells = @("C5", "C16", "C19", "C27", "C28")
foreach ($addr in $dvCells) {
    $dv = $rs0003.Range($addr).Validation
    $dv.InCellDropdown = $true
    $dv.ShowInput = $false
    $dv.ShowError = $false
}

# ---------------------------------------------------------------------------
# 3. performance_map: add operation_state lookup column (column E)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("performance_map")

# Header row (row 2) - matches D2's empty "lookup_variables" group style
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Variable name row (row 3) - matches D3's style, new lookup variable name
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E3").Value = "operation_state"
$ws.Range("E3").AddComment("The operation state at the operating conditions")

# Units row (row 4) - matches D4's style, unitless value
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E4").Value = "-"

# Data rows (5-144) - matches D5's style, every operating point is NORMAL
$ws.Range("D5").Copy()
$ws.Range("E5:E144").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
for ($r = 5; $r -le 144; $r++) {
    $ws.Cells.Item($r, 5).Value = "NORMAL"
}

Write-Output "done"
